$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ECs" target-cluster row (original row 2). This shifts the
# remaining rows (FAPs, MuSCs) up by one and also drops the now-unused
# "ECs" shared string from the workbook.
$ws.Rows(2).Delete()

# Recalculate the derived-specificity columns (O, P, S, T) now that the
# "ECs" row no longer contributes to the per-column totals.
$ws.Range("O2").Value = 0.9851515664921635
$ws.Range("P2").Value = 0.9851515664921635
$ws.Range("S2").Value = 0.9851515664921635
$ws.Range("T2").Value = 0.9851515664921635

$ws.Range("O3").Value = 0.01484843350783645
$ws.Range("P3").Value = 0.01484843350783645
$ws.Range("S3").Value = 0.01484843350783645
$ws.Range("T3").Value = 0.01484843350783645
